# Atualização dos gráficos 26082020
$wb = $excel.ActiveWorkbook

# --- Sheet "Mensal" update: row 14 ---
$wsMensal = $wb.Worksheets.Item("Mensal")
$wsMensal.Range("A14").Value = 44066
$wsMensal.Range("B14").Value = 131.55
$wsMensal.Range("D14").Value = -7.83

# --- Sheet "Diario" append rows 384-390 ---
$wsDiario = $wb.Worksheets.Item("Diario")

$wsDiario.Range("A384").Value = 44060
$wsDiario.Range("B384").Value = 199.4
$wsDiario.Range("C384").Value = 142.73
$wsDiario.Range("D384").Value = 39.7

$wsDiario.Range("A385").Value = 44061
$wsDiario.Range("B385").Value = 218.18
$wsDiario.Range("C385").Value = 142.73
$wsDiario.Range("D385").Value = 52.86

$wsDiario.Range("A386").Value = 44062
$wsDiario.Range("B386").Value = 236.66
$wsDiario.Range("C386").Value = 142.73
$wsDiario.Range("D386").Value = 65.81

$wsDiario.Range("A387").Value = 44063
$wsDiario.Range("B387").Value = 256.02
$wsDiario.Range("C387").Value = 142.73
$wsDiario.Range("D387").Value = 79.37

$wsDiario.Range("A388").Value = 44064
$wsDiario.Range("B388").Value = 276.91
$wsDiario.Range("C388").Value = 142.73
$wsDiario.Range("D388").Value = 94.01000000000001

$wsDiario.Range("A389").Value = 44065
$wsDiario.Range("B389").Value = 300.08
$wsDiario.Range("C389").Value = 142.73
$wsDiario.Range("D389").Value = 110.25

$wsDiario.Range("A390").Value = 44066
$wsDiario.Range("B390").Value = 313.94
$wsDiario.Range("C390").Value = 142.73
$wsDiario.Range("D390").Value = 119.95

# Copy date style (column A) from row 383 down to the new rows to match formatting (s="2")
$wsDiario.Range("A383").Copy() | Out-Null
$wsDiario.Range("A384:A390").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false
